$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.674.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.774.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '356.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.18'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.46%  '
$ws.Range('E7').Value = '  -3.17%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.584'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.52'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.89%  '
$ws.Range('E11').Value = '  +3.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0842'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.45'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.206.12'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.766.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.929'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.563.42'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.94'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  -2.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.37%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.165'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +14.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '51.42'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.40'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0447'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0837'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.16'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('E39').Value = '  -4.69%  '
$ws.Range('E40').Value = '  -4.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.54'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.31%  '
$ws.Range('E42').Value = '  -3.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.81%  '
$ws.Range('E45').Value = '  -6.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.086.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.938'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.55'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.54'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.52%  '
